$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2089.5881
$ws.Range("J70").Value = 2501.9167
$ws.Range("L70").Value = 7505.750100000001
$ws.Range("N70").Value = -8045.750100000001
$ws.Range("H73").Value = 2089.5881
$ws.Range("J73").Value = 2501.9167
$ws.Range("L73").Value = 7505.750100000001
$ws.Range("N73").Value = -9377.750100000001
$ws.Range("H87").Value = 41450
$ws.Range("J87").Value = 41450
$ws.Range("L87").Value = 41450
$ws.Range("N87").Value = -43946
$ws.Range("H90").Value = 41450
$ws.Range("J90").Value = 41450
$ws.Range("L90").Value = 124350
$ws.Range("N90").Value = -136830
$ws.Range("H100").Value = 33436174
$ws.Range("I100").Value = 41794092
$ws.Range("J100").Value = 4500
$ws.Range("K100").Value = 41794092
$ws.Range("L100").Value = 4500
$ws.Range("M100").Value = -41793551
$ws.Range("N100").Value = -5582
$ws.Range("H134").Value = 59833.332
$ws.Range("J134").Value = 59833.332
$ws.Range("L134").Value = 59833.332
$ws.Range("N134").Value = -69973.33199999999
$ws.Range("H139").Value = 20000
$ws.Range("I139").Value = 20000
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 20000
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -14860
$ws.Range("N139").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23807.389
$ws.Range("I32").Value = 3479.7446
$ws.Range("J32").Value = 501507
$ws.Range("K32").Value = 3479.7446
$ws.Range("L32").Value = 501507
$ws.Range("M32").Value = -3192.7446
$ws.Range("N32").Value = -502081
$ws.Range("H88").Value = 5320.3
$ws.Range("I88").Value = 2861.2
$ws.Range("J88").Value = 7779.4
$ws.Range("K88").Value = 2861.2
$ws.Range("L88").Value = 7779.4
$ws.Range("M88").Value = -2455.2
$ws.Range("N88").Value = -8591.4
$ws.Range("H91").Value = 5320.3
$ws.Range("I91").Value = 2861.2
$ws.Range("J91").Value = 7779.4
$ws.Range("K91").Value = 2861.2
$ws.Range("L91").Value = 7779.4
$ws.Range("M91").Value = -1457.2
$ws.Range("N91").Value = -10587.4

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H99").Value = 1573.7084
$ws.Range("I99").Value = 1494.0454
$ws.Range("K99").Value = 1494.0454
$ws.Range("M99").Value = 3.954600000000028

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2147
$ws.Range("I31").Value = 1280.5555
$ws.Range("K31").Value = 1280.5555
$ws.Range("M31").Value = -985.5554999999999
$ws.Range("H34").Value = 2147
$ws.Range("I34").Value = 1280.5555
$ws.Range("K34").Value = 1280.5555
$ws.Range("M34").Value = -1078.5555
$ws.Range("H58").Value = 3888.4285
$ws.Range("I58").Value = 2466.6667
$ws.Range("J58").Value = 4954.75
$ws.Range("K58").Value = 2466.6667
$ws.Range("L58").Value = 4954.75
$ws.Range("M58").Value = -2263.6667
$ws.Range("N58").Value = -5360.75
$ws.Range("H132").Value = 2234.0293
$ws.Range("I132").Value = 1887.4
$ws.Range("K132").Value = 5662.200000000001
$ws.Range("M132").Value = -3132.200000000001
$ws.Range("H134").Value = 1729.2532
$ws.Range("I134").Value = 1166.0526
$ws.Range("K134").Value = 3498.1578
$ws.Range("M134").Value = -963.1578
$ws.Range("H136").Value = 3888.4285
$ws.Range("I136").Value = 2466.6667
$ws.Range("J136").Value = 4954.75
$ws.Range("K136").Value = 7400.000100000001
$ws.Range("L136").Value = 14864.25
$ws.Range("M136").Value = -4850.000100000001
$ws.Range("N136").Value = -19964.25

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 8500
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H132").Value = 2667.3333
$ws.Range("I132").Value = 2502.238
$ws.Range("J132").Value = 3129.6
$ws.Range("K132").Value = 7506.714
$ws.Range("L132").Value = 9388.799999999999
$ws.Range("M132").Value = -4976.714
$ws.Range("N132").Value = -14448.8
$ws.Range("H134").Value = 19495.5
$ws.Range("J134").Value = 19495.5
$ws.Range("L134").Value = 58486.5
$ws.Range("N134").Value = -63556.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5561.6
$ws.Range("I61").Value = 5809.222
$ws.Range("J61").Value = 3333
$ws.Range("K61").Value = 5809.222
$ws.Range("L61").Value = 3333
$ws.Range("M61").Value = -5607.222
$ws.Range("N61").Value = -3737
$ws.Range("H113").Value = 5561.6
$ws.Range("I113").Value = 5809.222
$ws.Range("J113").Value = 3333
$ws.Range("K113").Value = 5809.222
$ws.Range("L113").Value = 3333
$ws.Range("M113").Value = -3639.222
$ws.Range("N113").Value = -7673
$ws.Range("H132").Value = 6842.4883
$ws.Range("I132").Value = 7512.7407
$ws.Range("K132").Value = 22538.2221
$ws.Range("M132").Value = -20008.2221
$ws.Range("H134").Value = 30989.8
$ws.Range("J134").Value = 30989.8
$ws.Range("L134").Value = 30989.8
$ws.Range("N134").Value = -41129.8
$ws.Range("H136").Value = 8995.056
$ws.Range("I136").Value = 7643.4287
$ws.Range("J136").Value = 9855.182000000001
$ws.Range("K136").Value = 22930.2861
$ws.Range("L136").Value = 29565.546
$ws.Range("M136").Value = -20380.2861
$ws.Range("N136").Value = -34665.546

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 12377
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 12377
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 12377
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -13157
$ws.Range("H117").Value = 40450
$ws.Range("J117").Value = 40450
$ws.Range("L117").Value = 40450
$ws.Range("N117").Value = -49628
$ws.Range("H132").Value = 11907563
$ws.Range("I132").Value = 21742186
$ws.Range("J132").Value = 2492.0527
$ws.Range("K132").Value = 65226558
$ws.Range("L132").Value = 7476.158100000001
$ws.Range("M132").Value = -65224028
$ws.Range("N132").Value = -12536.1581
